$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update ESTADO column values for the rows that received user feedback
$ws.Range("G8").Value = 1
$ws.Range("G10").Value = 1
$ws.Range("G11").Value = 1

# Move the active selection to G11 (keeps the frozen pane at row 2)
$ws.Range("G11").Select()

# Zoom the sheet to 70%
$excel.ActiveWindow.Zoom = 70
